$d = $word.ActiveDocument

function Replace-Text($old, $new) {
    $found = $d.Content.Find.Execute($old, $false, $false, $false, $false, $false, $true, 1, $false, $new, 2)
    if (-not $found) {
        Write-Output "NOT FOUND: $old"
    }
}

# 1. "Our guiding question is:" -> "Our guiding question:"
Replace-Text "Our guiding question is:" "Our guiding question:"

# 2. marine reserves sentence gains "and single speices closures"
Replace-Text "various sizes of marine reserves. The output" "various sizes of marine reserves and single speices closures. The output"

# 3. "We hope this model..." -> "This model..." + "multiple" -> "multi-species"
Replace-Text "We hope this model will be beneficial to resource mangaers and planners who use marine reserves with the main goal of multiple fisheries recovery." "This model will be beneficial to resource mangaers and planners who use marine reserves with the main goal of multi-species fisheries recovery."

# 4. Insert new "1. Predator-prey model" Heading5 before the "Submodels will include:" paragraph,
#    and collapse that paragraph's three runs into the simplified single-run text.
$target = $null
foreach ($p in $d.Paragraphs) {
    if ($p.Range.Text -like "Submodels will include:*") {
        $target = $p
        break
    }
}
if ($target -eq $null) {
    Write-Output "NOT FOUND: Submodels paragraph"
} else {
    $prevEnd = $target.Range.Start - 1
    $insertionPoint = $d.Range($prevEnd, $prevEnd)
    $headingFrag = @'
<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:pPr><w:pStyle w:val="Heading5"/></w:pPr><w:bookmarkStart w:id="24" w:name="predator-prey-model"/><w:bookmarkEnd w:id="24"/><w:r><w:t xml:space="preserve">1. Predator-prey model</w:t></w:r></w:p>
'@
    [void]$insertionPoint.InsertXML($headingFrag)
}

Replace-Text "Submodels will include: 1. Predator-prey model: To simulate interaction between predator and prey, including growth, harvest and competition among both species. This model will be lumped, dynamic, stochastic and abstract. Inputs will include biological parameters such as intrinsic growth rate, carrying capcity and competition coefficients, as well as, harvest rates. The model will utliize generalist predator-prey variables, thus sensitivity analysis will be conducted for intrinsic growth rates." "To simulate interaction between predator and prey, including growth, harvest and competition among species. This model will be lumped, dynamic, stochastic and abstract. Inputs will include biological parameters such as intrinsic growth rate, carrying capcity and competition coefficients, as well as, harvest rates (Table 1). The model will utliize generalist predator-prey variables, thus sensitivity analysis will be conducted for intrinsic growth rates."

# 5. "...mathematicaly described as follows:" -> "...follows (Samhouri et al., 2017):"
Replace-Text "Interactions of generalist predator prey dynamics are mathematicaly described as follows:" "Interactions of generalist predator prey dynamics are mathematicaly described as follows (Samhouri et al., 2017):"

# 6. Delete the entire "Here, X and P denote..." explanatory paragraph (with embedded oMath runs).
$target = $null
foreach ($p in $d.Paragraphs) {
    if ($p.Range.Text -like "Here, X and P denote*") {
        $target = $p
        break
    }
}
if ($target -eq $null) {
    Write-Output "NOT FOUND: Here X and P paragraph"
} else {
    $target.Range.Delete()
}

Write-Output "done phase 3"

# 7. Replace the old flat numbered-list sub-model descriptions with the new
#    Heading5-structured sections (Marine reserve model, Economic Model, Wrapper).
$first = $null
$last = $null
foreach ($p in $d.Paragraphs) {
    if ($p.Range.Text -like "MPA model*") {
        $first = $p
    }
    if ($p.Range.Text -like "Determine the inputs*") {
        $last = $p
    }
}
if (($first -eq $null) -or ($last -eq $null)) {
    Write-Output "NOT FOUND: list block"
} else {
    $full = $d.Range($first.Range.Start, $last.Range.End)
    $listFrag = @'
<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main" xmlns:m="http://schemas.openxmlformats.org/officeDocument/2006/math"><w:pPr><w:pStyle w:val="Heading5"/></w:pPr><w:bookmarkStart w:id="25" w:name="marine-reserve-model"/><w:bookmarkEnd w:id="25"/><w:r><w:t xml:space="preserve">2. Marine reserve model</w:t></w:r></w:p><w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:pPr><w:pStyle w:val="FirstParagraph"/></w:pPr><w:r><w:t xml:space="preserve">A patch model to simulate predator-prey interactions given spatial closures. Space will be represented by a series of vectors. Sensitivity analysis will be consucted on the size of the closure, represented by number of patches closed where harvest of prey and predator = 0. The model will be spatial, static, stocastic and abstract.</w:t></w:r></w:p><w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:pPr><w:pStyle w:val="Heading5"/></w:pPr><w:bookmarkStart w:id="26" w:name="economic-model"/><w:bookmarkEnd w:id="26"/><w:r><w:t xml:space="preserve">3. Economic Model</w:t></w:r></w:p><w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:pPr><w:pStyle w:val="FirstParagraph"/></w:pPr><w:r><w:t xml:space="preserve">To generate net present value generated by harvest of predator and prey. Inputs to this model will include price of species ($/individual), amount harvested (# of individuals), disocunt rate (%) and time (years). The model will be determinitic.</w:t></w:r></w:p><w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:pPr><w:pStyle w:val="BodyText"/></w:pPr><w:r><w:t xml:space="preserve">Net present value (NPV) of harvest from predator and prey will be calcualted using the following equation:</w:t></w:r></w:p><w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main" xmlns:m="http://schemas.openxmlformats.org/officeDocument/2006/math"><w:pPr><w:pStyle w:val="BodyText"/></w:pPr><m:oMathPara><m:oMathParaPr><m:jc m:val="center"/></m:oMathParaPr><m:oMath><m:r><m:t>N</m:t></m:r><m:r><m:t>P</m:t></m:r><m:r><m:t>V</m:t></m:r><m:r><m:t>=</m:t></m:r><m:f><m:fPr><m:type m:val="bar"/></m:fPr><m:num><m:r><m:t>P</m:t></m:r><m:r><m:t>H</m:t></m:r></m:num><m:den><m:r><m:t>(</m:t></m:r><m:r><m:t>1</m:t></m:r><m:r><m:t>+</m:t></m:r><m:r><m:t>D</m:t></m:r><m:sSup><m:e><m:r><m:t>)</m:t></m:r></m:e><m:sup><m:r><m:t>t</m:t></m:r></m:sup></m:sSup></m:den></m:f></m:oMath></m:oMathPara></w:p><w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:pPr><w:pStyle w:val="FirstParagraph"/></w:pPr><w:r><w:t xml:space="preserve">Where P is the market price, H is harvest, D is the discount rate and t represent time in years.</w:t></w:r></w:p><w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:pPr><w:pStyle w:val="Heading5"/></w:pPr><w:bookmarkStart w:id="27" w:name="wrapper-of-these-three-models"/><w:bookmarkEnd w:id="27"/><w:r><w:t xml:space="preserve">4. Wrapper of these three models</w:t></w:r></w:p><w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:pPr><w:pStyle w:val="FirstParagraph"/></w:pPr><w:r><w:t xml:space="preserve">A wrapper function will be used to simulate predator-prey dynamics given the marine reserve model implementation and calculate economic returns. The function will also output relevant graphs over projected time.</w:t></w:r></w:p><w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:pPr><w:pStyle w:val="BodyText"/></w:pPr><w:r><w:t xml:space="preserve">Model simulations will follow the following steps:</w:t></w:r><w:r><w:br w:type="textWrapping"/></w:r><w:r><w:t xml:space="preserve">1. Allow the system to run until equilibrum with no intervention</w:t></w:r><w:r><w:br w:type="textWrapping"/></w:r><w:r><w:t xml:space="preserve">2. Siulate trophic downgrading by harvesting the predator until equilibrium and then harvesting the prey until a new equilibrium is reached.</w:t></w:r><w:r><w:br w:type="textWrapping"/></w:r><w:r><w:t xml:space="preserve">3. Implement one of the following management strategies:</w:t></w:r><w:r><w:br w:type="textWrapping"/></w:r><w:r><w:t xml:space="preserve">+i) Harvest of both species is reduced due to marine reserve implementation at various sizes, where a marine reserve of 100% of the area results in harvesting of both species be equal to 0 and thus stopped completely.</w:t></w:r><w:r><w:br w:type="textWrapping"/></w:r><w:r><w:t xml:space="preserve">+ii) Harvest of predator stops (i.e. hp = 0) until new eqilibrium is reached, and then harvest of prey stops (i.e. hx = 0)</w:t></w:r><w:r><w:br w:type="textWrapping"/></w:r><w:r><w:t xml:space="preserve">+iii)Harvest of prey stops (i.e. hx = 0) until new eqilibrium is reached, and then harvest of predator stops (i.e. hp = 0)</w:t></w:r></w:p>
'@
    [void]$full.InsertXML($listFrag)
}

Write-Output "done phase 4"
